$d = $word.ActiveDocument

# 1. "...كيفية عملها." -> "...كيفية عمل ParentText."
$d.Content.Find.Execute(
    "قبل أن تتلقوا النصائح، خلونا نلقي نَّظْرَة على كيفية عملها. ", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "قبل أن تتلقوا النصائح، خلونا نلقي نَّظْرَة على كيفية عمل ParentText. ", 2) | Out-Null

# 2. "غير متحقِّق" -> "غير متأكد" (word-level match to avoid AutoCorrect touching the straight quotes later in the sentence)
$d.Content.Find.Execute(
    "غير متحقِّق من", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "غير متأكد من", 2) | Out-Null

# 3. "لمشاركة رابط هذا الروبوت..." -> "لمشاركة رابط روبوت الدردشة... ParentText"
$d.Content.Find.Execute(
    "لمشاركة رابط هذا الروبوت مع صديق، اختر “ادعُ صديقًا ل", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "لمشاركة رابط روبوت الدردشة مع صديق، اختر “ادعُ صديقًا ل ParentText”", 2) | Out-Null

# 4. "ادعُ صديق إلى برنامج رسائل الأهل" -> "ادعُ صديق إلى برنامج ParentText"
$d.Content.Find.Execute(
    "ادعُ صديق إلى برنامج رسائل الأهل", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "ادعُ صديق إلى برنامج ParentText", 2) | Out-Null

# 5. "شاهد فيديو عن برنامج رسائل الأهل" -> "شاهد فيديو عن برنامج ParentText"
$d.Content.Find.Execute(
    "شاهد فيديو عن برنامج رسائل الأهل", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "شاهد فيديو عن برنامج ParentText", 2) | Out-Null
